# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (fund-holdings detail) right after
#    "2021-Q4" and before "总计".
# 2. Refresh the "总计" (summary) worksheet with a new leading row for
#    2022-Q1 (existing rows shift down by one).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# Step 0: remove the existing "总计" sheet. Because the engine hands out
# sheetId values from a monotonically increasing counter that currently
# sits at 4 (the workbook has exactly 4 sheets, ids 1-4), deleting "总计"
# (id 4) and immediately minting the next two sheets reproduces the
# target numbering: "2022-Q1" -> sheetId 4, "总计" -> sheetId 5.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Delete()

# ---------------------------------------------------------------------------
# Step 1: create "2022-Q1" by duplicating "2021-Q4" (same column layout:
# 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名,
# same header/index-column styling) right after "2021-Q4", then drop the
# old sheet's extra data rows and overwrite rows 2-33 with the new values.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newFund = $wb.Worksheets.Item(4)
$newFund.Name = "2022-Q1"

# "2021-Q4" has 59 data rows (r2:r60); the new sheet only needs 32 (r2:r33).
$newFund.Rows("34:60").Delete()

$fundData = @(
    @("110005","易方达积极成长混合","36.95","92.57","7.23","2.6715",5),
    @("506005","博时科创板三年定期开放混合","22.84","96.44","3.74","0.8542",7),
    @("009663","华泰紫金科技创新3年封闭运作灵活配置混合A","10.78","68.04","4.59","0.4948",4),
    @("501202","华泰紫金科技创新3年封闭运作灵活配置混合C","10.01","68.04","4.59","0.4595",4),
    @("213003","宝盈策略增长混合","10.28","94.38","4.33","0.4451",10),
    @("213002","宝盈泛沿海增长混合","5.10","93.76","6.17","0.3147",6),
    @("008866","博时产业新趋势灵活配置混合A","7.89","86.47","3.37","0.2659",7),
    @("501201","红土创新科技创新 3 年封闭运作灵活配置混合","3.99","96.70","4.07","0.1624",4),
    @("008811","鹏华科技创新混合","3.49","94.05","3.76","0.1312",6),
    @("501098","建信科技创新 3 年封闭运作灵活配置混合","3.61","76.11","2.97","0.1072",7),
    @("519967","长信利富债券","5.83","20.15","1.21","0.0705",1),
    @("580006","东吴新经济混合","1.17","91.56","5.52","0.0646",1),
    @("673040","西部利得行业主题优选灵活配置混合A","4.73","29.77","0.78","0.0369",6),
    @("200016","长城稳健成长灵活配置混合","0.83","78.39","4.22","0.0350",6),
    @("690003","民生加银精选混合","0.61","91.83","5.72","0.0349",4),
    @("168401","红土创新转型精选灵活配置混合（LOF）","0.78","93.82","3.99","0.0311",5),
    @("000994","建信睿盈灵活配置混合A","0.58","88.02","5.06","0.0293",4),
    @("163818","中银中小盘成长混合","0.98","87.49","2.64","0.0259",4),
    @("004332","恒生前海沪港深新兴产业精选混合","0.52","80.98","4.47","0.0232",5),
    @("673043","西部利得行业主题优选灵活配置混合C","2.67","29.77","0.78","0.0208",6),
    @("008082","国寿安保研究精选混合A","0.52","91.60","3.65","0.0190",8),
    @("007965","民生加银品质消费股票A","0.25","88.39","5.47","0.0137",5),
    @("008867","博时产业新趋势灵活配置混合C","0.40","86.47","3.37","0.0135",7),
    @("005437","易方达易百智能量化策略灵活配置混合A","0.96","94.52","1.19","0.0114",2),
    @("710002","富安达策略精选混合","0.63","63.45","1.58","0.0100",6),
    @("000995","建信睿盈灵活配置混合C","0.19","88.02","5.06","0.0096",4),
    @("013371","民生加银新能源智选混合A","0.37","78.61","2.14","0.0079",10),
    @("007966","民生加银品质消费股票C","0.13","88.39","5.47","0.0071",5),
    @("740001","长安宏观策略混合","0.16","71.93","4.14","0.0066",7),
    @("008083","国寿安保研究精选混合C","0.15","91.60","3.65","0.0055",8),
    @("013372","民生加银新能源智选混合C","0.13","78.61","2.14","0.0028",10),
    @("005438","易方达易百智能量化策略灵活配置混合C","0.17","94.52","1.19","0.0020",2)
)

# Columns B-G (fund code, fund name, fund size, total equity position,
# position ratio, held market value) are stored as text in the source
# workbook even when they look numeric - force text storage so Excel
# doesn't silently coerce them (and lose precision/leading context).
$newFund.Range("B2:G33").NumberFormat = "@"

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $r = $i + 2
    $row = $fundData[$i]
    $newFund.Cells.Item($r, 2).Value = $row[0]   # B: 基金代码
    $newFund.Cells.Item($r, 3).Value = $row[1]   # C: 基金名称
    $newFund.Cells.Item($r, 4).Value = $row[2]   # D: 基金规模
    $newFund.Cells.Item($r, 5).Value = $row[3]   # E: 股票总仓位
    $newFund.Cells.Item($r, 6).Value = $row[4]   # F: 仓位占比
    $newFund.Cells.Item($r, 7).Value = $row[5]   # G: 持有市值(亿元)
    $newFund.Cells.Item($r, 8).Value = [int]$row[6]  # H: 仓位排名 (number)
}

# ---------------------------------------------------------------------------
# Step 2: recreate "总计" at the end of the tab strip (claims sheetId 5),
# matching the style used by the other summary-style sheets (index-column
# + header formatting borrowed from "2021-Q3", which is untouched and
# still carries that exact style).
# ---------------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add($null, $newFund)
$newTotal.Name = "总计"

$fmtSrc = $wb.Worksheets.Item("2021-Q3")
$fmtSrc.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$fmtSrc.Range("A2:A5").Copy()
$newTotal.Range("A2:A5").PasteSpecial(-4122)

$newTotal.Cells.Item(1, 2).Value = "日期"
$newTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$newTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 32, 6.39),
    @(1, "2021-Q4", 59, 20.52),
    @(2, "2021-Q3", 22, 9.18),
    @(3, "2021-Q2", 6, 0.35)
)
for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $newTotal.Cells.Item($r, 1).Value = [int]$row[0]
    $newTotal.Cells.Item($r, 2).Value = $row[1]
    $newTotal.Cells.Item($r, 3).Value = [int]$row[2]
    $newTotal.Cells.Item($r, 4).Value = $row[3]
}

Write-Output "2022-Q1 + 总计 rebuilt"
